$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = '0.984 (0.984 ± 0.000)'
$ws.Range('C2').Value = '00:04:50 (00:05:02 ± 00:00:08)'
$ws.Range('D2').Value = '00:00:04 (00:00:06 ± 00:00:01)'
$ws.Range('B3').Value = '0.991 (0.982 ± 0.005)'
$ws.Range('C3').Value = '00:00:25 (00:01:27 ± 00:00:47)'
$ws.Range('D3').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B4').Value = '0.987 (0.974 ± 0.007)'
$ws.Range('C4').Value = '00:00:46 (00:01:04 ± 00:00:12)'
$ws.Range('D4').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B5').Value = '0.987 (0.973 ± 0.009)'
$ws.Range('C5').Value = '00:05:11 (00:05:22 ± 00:00:14)'
$ws.Range('D5').Value = '00:00:00 (00:00:01 ± 00:00:00)'
$ws.Range('B6').Value = '0.993 (0.985 ± 0.006)'
$ws.Range('C6').Value = '00:04:56 (00:05:04 ± 00:00:04)'
$ws.Range('D6').Value = '00:00:00 (00:00:05 ± 00:00:02)'
$ws.Range('B7').Value = '0.989 (0.976 ± 0.006)'
$ws.Range('C7').Value = '00:05:00 (00:05:02 ± 00:00:02)'
$ws.Range('D7').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B8').Value = '0.987 (0.980 ± 0.004)'
$ws.Range('C8').Value = '00:04:22 (00:06:07 ± 00:01:49)'
$ws.Range('D8').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B9').Value = '0.994 (0.984 ± 0.005)'
$ws.Range('C9').Value = '00:04:59 (00:05:01 ± 00:00:02)'
$ws.Range('D9').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B10').Value = '0.991 (0.981 ± 0.007)'
$ws.Range('C10').Value = '00:04:29 (00:04:29 ± 00:00:00)'
$ws.Range('D10').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B11').Value = '0.991 (0.983 ± 0.005)'
$ws.Range('C11').Value = '00:05:05 (00:05:06 ± 00:00:00)'
$ws.Range('D11').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B12').Value = '0.987 (0.977 ± 0.006)'
$ws.Range('C12').Value = '00:02:10 (00:02:34 ± 00:00:17)'
$ws.Range('D12').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B13').Value = '0.991 (0.976 ± 0.008)'
$ws.Range('C13').Value = '00:00:02 (00:00:03 ± 00:00:00)'
$ws.Range('D13').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B14').Value = '0.991 (0.979 ± 0.006)'
$ws.Range('C14').Value = '00:00:41 (00:00:45 ± 00:00:01)'
$ws.Range('D14').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B15').Value = '0.987 (0.981 ± 0.006)'
$ws.Range('C15').Value = '00:00:01 (00:04:11 ± 00:01:42)'
$ws.Range('D15').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B16').Value = '0.994 (0.983 ± 0.005)'
$ws.Range('C16').Value = '00:01:05 (00:01:12 ± 00:00:04)'
$ws.Range('D16').Value = '00:00:00 (00:00:00 ± 00:00:00)'
$ws.Range('B17').Value = '0.987 (0.970 ± 0.017)'
$ws.Range('C17').Value = '00:00:30 (00:01:57 ± 00:01:31)'
$ws.Range('D17').Value = '00:00:00 (00:00:00 ± 00:00:00)'
